$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

$data = @{
    2  = @(2196.65, 2149.85)
    3  = @(424.6, 425.85)
    4  = @(1610.4, 1593.8)
    5  = @(7103.15, 7068.6)
    6  = @(213.9, 213.95)
    7  = @(195.75, 197.65)
    8  = @(43791.85, 43742.15)
    9  = @(634.45, 626.25)
    10 = @(3844.2, 3846.45)
    11 = @(146.7, 147.4)
    12 = @(1326.6, 1325.4)
    13 = @(506.85, 497.3)
    14 = @(1502.05, 1504.8)
    15 = @(654.65, 649.8)
    16 = @(462.05, 464.6)
    17 = @(1559.75, 1552.2)
    18 = @(275.15, 275.45)
    19 = @(19841.3, 19751.2)
    20 = @(250, 251.05)
    21 = @(561.5, 563.75)
    22 = @(662.55, 680.15)
    23 = @(681.7, 674.3)
    24 = @(262.65, 260.85)
    25 = @(126.25, 124.35)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}
